$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update title / label text (October 2016 -> November 2016) ---
$ws.Range("A2").Value2 = "by Sector, 2006-November 2016 (Thousand Tons)"

# --- Insert a new row for "November" month data right before the       ---
# --- existing "Annual Totals" row (currently row 53), shifting         ---
# --- everything below it down by one row.                              ---
$ws.Rows(53).Insert()

# Copy number formatting from the row above (October data row, still 52)
# onto the newly inserted row so it matches the rest of the monthly rows
# (style only - values are overwritten right after).
$ws.Range("A52:F52").Copy()
$ws.Range("A53:F53").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A53").Value2 = "November"
$ws.Range("B53").Value2 = 1294
$ws.Range("C53").Value2 = 39
$ws.Range("D53").Value2 = 1069
$ws.Range("E53").Value2 = 184
$ws.Range("F53").Value2 = 0.43

# --- Update the "Year 2014" / "Year 2015" / "Year 2016" annual rows    ---
# --- under "Annual Totals" (now rows 55-57 after the insert above).    ---
$ws.Range("B55").Value2 = 15309
$ws.Range("C55").Value2 = 408
$ws.Range("D55").Value2 = 12654
$ws.Range("E55").Value2 = 2242
$ws.Range("F55").Value2 = 5

$ws.Range("B56").Value2 = 15160
$ws.Range("C56").Value2 = 415
$ws.Range("D56").Value2 = 12565
$ws.Range("E56").Value2 = 2173
$ws.Range("F56").Value2 = 7

$ws.Range("B57").Value2 = 14658
$ws.Range("C57").Value2 = 425
$ws.Range("D57").Value2 = 12176
$ws.Range("E57").Value2 = 2050
$ws.Range("F57").Value2 = 7

# --- Update the "Rolling 12 Months Ending in ..." label (now row 58). ---
$ws.Range("A58").Value2 = "Rolling 12 Months Ending in November"

# --- Update the "Year to Date" annual rows (now rows 59-60).           ---
$ws.Range("B59").Value2 = 16556
$ws.Range("C59").Value2 = 451
$ws.Range("D59").Value2 = 13720
$ws.Range("E59").Value2 = 2378
$ws.Range("F59").Value2 = 7

$ws.Range("B60").Value2 = 16129
$ws.Range("C60").Value2 = 462
$ws.Range("D60").Value2 = 13408
$ws.Range("E60").Value2 = 2252
$ws.Range("F60").Value2 = 8
